$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while forcing Excel to store it as
# TEXT (matches the source data, which is plain inlineStr/text, not
# numeric) - uses the classic leading-apostrophe "quote prefix".
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

Set-TextCell "D2" '25.837.63'
$ws.Range("E2").Value = '  -0.41%  '
Set-TextCell "D3" '1.636.96'
$ws.Range("E3").Value = '  -0.15%  '
Set-TextCell "D4" '1.003'
$ws.Range("E4").Value = '  +0.13%  '
Set-TextCell "D5" '215.86'
Set-TextCell "D6" '0.5065'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.14%  '
Set-TextCell "D8" '0.2576'
$ws.Range("E8").Value = '  +0.08%  '
Set-TextCell "D9" '0.06429'
$ws.Range("E9").Value = '  +1.24%  '
$ws.Range("E10").Value = '  -1.30%  '
Set-TextCell "D11" '0.07781'
$ws.Range("E11").Value = '  +0.65%  '
Set-TextCell "D12" '4.284'
$ws.Range("E12").Value = '  -0.30%  '
Set-TextCell "D13" '1.862.55'
$ws.Range("E13").Value = '  -0.16%  '
Set-TextCell "D14" '1.638.01'
$ws.Range("E14").Value = '  -0.46%  '
Set-TextCell "D15" '0.5618'
$ws.Range("E15").Value = '  +2.92%  '
Set-TextCell "D16" '63.21'
$ws.Range("E16").Value = '  -1.59%  '
Set-TextCell "D17" '0.0₅7588'
$ws.Range("E17").Value = '  -1.93%  '
Set-TextCell "D18" '25.853.19'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  +0.14%  '
Set-TextCell "D20" '194.39'
$ws.Range("E20").Value = '  -0.57%  '
Set-TextCell "D21" '4.330'
$ws.Range("E21").Value = '  -3.02%  '
Set-TextCell "D22" '9.868'
$ws.Range("E22").Value = '  -0.82%  '
Set-TextCell "D23" '6.082'
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("E24").Value = '  +0.02%  '
Set-TextCell "D25" '1.790'
$ws.Range("E25").Value = '  -5.79%  '
Set-TextCell "D26" '0.1275'
$ws.Range("E26").Value = '  +0.55%  '
Set-TextCell "D27" '140.11'
$ws.Range("E27").Value = '  -1.85%  '
Set-TextCell "D28" '6.775'
$ws.Range("E28").Value = '  -1.13%  '
Set-TextCell "D29" '15.45'
$ws.Range("E29").Value = '  -0.73%  '
Set-TextCell "D30" '1.240'
$ws.Range("E30").Value = '  +0.22%  '
Set-TextCell "D31" '0.04879'
$ws.Range("E31").Value = '  -0.18%  '
Set-TextCell "D32" '3.295'
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  +0.39%  '
Set-TextCell "D34" '1.555'
$ws.Range("E34").Value = '  +0.20%  '
Set-TextCell "D35" '2.365'
$ws.Range("E35").Value = '  -0.41%  '
Set-TextCell "D36" '0.9019'
$ws.Range("E36").Value = '  -1.67%  '
Set-TextCell "D37" '2.575'
$ws.Range("E37").Value = '  +0.28%  '
Set-TextCell "D38" '1.129.98'
$ws.Range("E38").Value = '  -0.40%  '
Set-TextCell "D39" '0.5494'
$ws.Range("E39").Value = '  -0.57%  '
Set-TextCell "D40" '0.01560'
$ws.Range("E40").Value = '  -0.39%  '
Set-TextCell "D41" '1.001'
$ws.Range("E41").Value = '  -0.04%  '
Set-TextCell "D42" '5.540'
$ws.Range("E42").Value = '  -0.78%  '
Set-TextCell "D43" '0.8004'
$ws.Range("E43").Value = '  -0.42%  '
Set-TextCell "D44" '97.44'
$ws.Range("E44").Value = '  -1.16%  '
Set-TextCell "D45" '1.789.53'
$ws.Range("E45").Value = '  +0.81%  '
Set-TextCell "D46" '0.0₈114'
$ws.Range("E46").Value = '  -6.16%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D47" '0.4445'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D48" '55.49'
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell "D49" '7.697'
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell "D50" '0.05052'
$ws.Range("E50").Value = '  -2.63%  '
Set-TextCell "D51" '1.005'
$ws.Range("E51").Value = '  +0.61%  '
